$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "ᶻᵍˣ"
$ws.Range("C28").Value = '"㊥Smok3y 1nOnly"'
$ws.Range("C42").Value = "㊥☆梅海听雪☆zgx"
$ws.Range("C64").Value = "囼uu文琪"
$ws.Range("C123").Value = "Globalking1001"
